# Oppdaterte omkodingstabell med dagens fylkesnummer
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before the old "rhf_kode" column (old C), shifting
# old C:E (rhf_kode, rhf_tekst, rhf_resh) to E:G. The former A:B
# (fylke_kode/fylke_tekst) stay in place and become fylker_historisk /
# fylker_tekst.
$ws.Range("C1:D1").EntireColumn.Insert()

# Rename the two pre-existing headers and set headers for the two newly
# inserted columns.
$ws.Cells.Item(1, 1).Value = "fylker_historisk"
$ws.Cells.Item(1, 2).Value = "fylker_tekst"
$ws.Cells.Item(1, 3).Value = "fylker_oppdatert"
$ws.Cells.Item(1, 4).Value = "ltmv"

# fylker_oppdatert (col C) / ltmv (col D) values per historic fylke (col A),
# keyed by the row's existing fylke code.
$oppdatert = @{
    1  = 1
    2  = 2
    3  = 3
    4  = 4
    5  = 5
    6  = 6
    7  = 7
    8  = 8
    9  = 9
    10 = 10
    11 = 11
    12 = 12
    14 = 14
    15 = 15
    16 = 50
    17 = 50
    18 = 18
    19 = 19
    20 = 20
    21 = 21
}
$ltmv = @{
    1  = 18
    2  = 17
    3  = 16
    4  = 15
    5  = 14
    6  = 13
    7  = 12
    8  = 11
    9  = 10
    10 = 9
    11 = 8
    12 = 7
    14 = 6
    15 = 5
    16 = 4
    17 = 4
    18 = 3
    19 = 2
    20 = 1
    21 = 1
}

for ($r = 2; $r -le 21; $r++) {
    $fylkeKode = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 3).Value = $oppdatert[[int]$fylkeKode]
    $ws.Cells.Item($r, 4).Value = $ltmv[[int]$fylkeKode]
}

# New row for the merged "Trøndelag" county (fylkesnummer 50), appended
# after the existing 20 data rows.
$ws.Cells.Item(22, 1).Value = 50
$ws.Cells.Item(22, 2).Value = "Trøndelag"
$ws.Cells.Item(22, 3).Value = 50
$ws.Cells.Item(22, 4).Value = 4
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = "Helse Midt-Norge"
$ws.Cells.Item(22, 7).Value = 100024

$ws.Range("D22").Select()
